$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of row 7 (A7:G7) while preserving the cell styles.
$ws.Range("A7:G7").ClearContents()
